{"js": "// Delete the \"League Center\" glossary entry paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(\"League Center:\") !== -1) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Delete the \"League Center\" glossary entry paragraph.\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -like \"*League Center:*\") {\n        $para.Range.Delete()\n    }\n}\n"}
